$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9091008
$ws.Range("I9").Value = 25000070
$ws.Range("J9").Value = 115.57143
$ws.Range("K9").Value = 25000070
$ws.Range("L9").Value = 115.57143
$ws.Range("M9").Value = -24999901
$ws.Range("N9").Value = -453.57143
$ws.Range("H33").Value = 115.85714
$ws.Range("I33").Value = 118.5
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 118.5
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 110.5
$ws.Range("N33").Value = -558
$ws.Range("H43").Value = 427.75
$ws.Range("I43").Value = 434.7143
$ws.Range("J43").Value = 422.33334
$ws.Range("K43").Value = 434.7143
$ws.Range("L43").Value = 422.33334
$ws.Range("M43").Value = -365.7143
$ws.Range("N43").Value = -560.33334
$ws.Range("H62").Value = 2910.7273
$ws.Range("I62").Value = 2923.6
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 2923.6
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -2299.6
$ws.Range("N62").Value = -4148
$ws.Range("H65").Value = 2910.7273
$ws.Range("I65").Value = 2923.6
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 14618
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -11498
$ws.Range("N65").Value = -20740
$ws.Range("H92").Value = 3401.25
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 6102.5
$ws.Range("K92").Value = 700
$ws.Range("L92").Value = 6102.5
$ws.Range("M92").Value = 548
$ws.Range("N92").Value = -8598.5
$ws.Range("H96").Value = 812.875
$ws.Range("I96").Value = 502.45456
$ws.Range("J96").Value = 1495.8
$ws.Range("K96").Value = 1507.36368
$ws.Range("L96").Value = 4487.4
$ws.Range("M96").Value = -134.3636799999999
$ws.Range("N96").Value = -7233.4
$ws.Range("H97").Value = 326.36365
$ws.Range("I97").Value = 266.66666
$ws.Range("J97").Value = 348.75
$ws.Range("K97").Value = 799.9999799999999
$ws.Range("L97").Value = 1046.25
$ws.Range("M97").Value = -303.9999799999999
$ws.Range("N97").Value = -2038.25
$ws.Range("H116").Value = 10992654
$ws.Range("I116").Value = 12823930
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 12823930
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -12820488
$ws.Range("N116").Value = -11884
$ws.Range("H125").Value = 750
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -2040
$ws.Range("N125").Value = -13920
$ws.Range("H137").Value = 146042
$ws.Range("I137").Value = 3058.2
$ws.Range("J137").Value = 503501.5
$ws.Range("K137").Value = 9174.599999999999
$ws.Range("L137").Value = 1510504.5
$ws.Range("M137").Value = -6624.599999999999
$ws.Range("N137").Value = -1515604.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4251.645
$ws.Range("I31").Value = 3570
$ws.Range("K31").Value = 3570
$ws.Range("M31").Value = -3275
$ws.Range("H34").Value = 4251.645
$ws.Range("I34").Value = 3570
$ws.Range("K34").Value = 3570
$ws.Range("M34").Value = -3368

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 444856
$ws.Range("I11").Value = 666734.2
$ws.Range("J11").Value = 1099.6666
$ws.Range("K11").Value = 2000202.6
$ws.Range("L11").Value = 3298.9998
$ws.Range("M11").Value = -2000062.6
$ws.Range("N11").Value = -3578.9998
$ws.Range("H12").Value = 69.2
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 74
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 222
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -568
$ws.Range("H38").Value = 67.75
$ws.Range("I38").Value = 67.75
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 203.25
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = 143.75
$ws.Range("H50").Value = 1027.5
$ws.Range("I50").Value = 55
$ws.Range("K50").Value = 165
$ws.Range("M50").Value = 316
$ws.Range("H52").Value = 716.5
$ws.Range("J52").Value = 716.5
$ws.Range("L52").Value = 2149.5
$ws.Range("N52").Value = -2681.5
$ws.Range("H53").Value = 1027.5
$ws.Range("I53").Value = 55
$ws.Range("K53").Value = 165
$ws.Range("M53").Value = 316
$ws.Range("H93").Value = 4892.316
$ws.Range("I93").Value = 3024
$ws.Range("J93").Value = 4996.1113
$ws.Range("K93").Value = 9072
$ws.Range("L93").Value = 14988.3339
$ws.Range("M93").Value = -7200
$ws.Range("N93").Value = -18732.3339
$ws.Range("H119").Value = 1351.125
$ws.Range("I119").Value = 1351.125
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 4053.375
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = 784.625
$ws.Range("H121").Value = 3334248
$ws.Range("I121").Value = 437.5
$ws.Range("J121").Value = 3847142
$ws.Range("K121").Value = 1312.5
$ws.Range("L121").Value = 11541426
$ws.Range("M121").Value = -2.5
$ws.Range("N121").Value = -11544046
$ws.Range("H137").Value = 7147.143
$ws.Range("I137").Value = 2353.3333
$ws.Range("J137").Value = 8454.546
$ws.Range("K137").Value = 7059.999899999999
$ws.Range("L137").Value = 25363.638
$ws.Range("M137").Value = -1959.999899999999
$ws.Range("N137").Value = -35563.638
$ws.Range("M38").ClearContents()
$ws.Range("M119").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4765062
$ws.Range("I122").Value = 6669120.5
$ws.Range("K122").Value = 20007361.5
$ws.Range("M122").Value = -20004911.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1692.6666
$ws.Range("I7").Value = 1191.375
$ws.Range("J7").Value = 2265.5715
$ws.Range("K7").Value = 1191.375
$ws.Range("L7").Value = 2265.5715
$ws.Range("M7").Value = -1079.375
$ws.Range("N7").Value = -2489.5715
$ws.Range("H22").Value = 715136
$ws.Range("I22").Value = 1429028.9
$ws.Range("J22").Value = 1243.1428
$ws.Range("K22").Value = 1429028.9
$ws.Range("L22").Value = 1243.1428
$ws.Range("M22").Value = -1428733.9
$ws.Range("N22").Value = -1833.1428
$ws.Range("H27").Value = 715136
$ws.Range("I27").Value = 1429028.9
$ws.Range("J27").Value = 1243.1428
$ws.Range("K27").Value = 1429028.9
$ws.Range("L27").Value = 1243.1428
$ws.Range("M27").Value = -1428921.9
$ws.Range("N27").Value = -1457.1428
$ws.Range("H40").Value = 4541.9165
$ws.Range("I40").Value = 2750.5
$ws.Range("K40").Value = 2750.5
$ws.Range("M40").Value = -2614.5
$ws.Range("H46").Value = 1134
$ws.Range("I46").Value = 927.6
$ws.Range("J46").Value = 1650
$ws.Range("K46").Value = 927.6
$ws.Range("L46").Value = 1650
$ws.Range("M46").Value = -739.6
$ws.Range("N46").Value = -2026
$ws.Range("H126").Value = 1692.6666
$ws.Range("I126").Value = 1191.375
$ws.Range("J126").Value = 2265.5715
$ws.Range("K126").Value = 3574.125
$ws.Range("L126").Value = 6796.7145
$ws.Range("M126").Value = -1104.125
$ws.Range("N126").Value = -11736.7145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1150.6364
$ws.Range("I126").Value = 873
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 2619
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -149
$ws.Range("N126").Value = -12140
